# Adds a "2022-Q1" worksheet (per-fund holdings) right before the "总计"
# (totals) sheet, and refreshes "总计" with a new leading row summarizing
# 2022-Q1 (existing rows shift down).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# Style-donor sheet: "2021-Q4" already uses the bold/centered/bordered
# header style (style index 2 in the original file) that every other
# per-fund sheet + the totals sheet share.
$styleDonor = $wb.Worksheets.Item("2021-Q4")

$totalOld = $wb.Worksheets.Item("总计")

# Remove the old "总计" sheet - we'll rebuild it fresh after inserting the
# new per-fund sheet so the sheetId numbering matches: 2022-Q1 gets the id
# freed up by the delete, 总计 gets a fresh (higher) id when re-added.
$totalOld.Delete() | Out-Null

# ---------------------------------------------------------------------
# New sheet: "2022-Q1" (per-fund holdings), placed at the end (i.e. right
# before where "总计" will be re-added).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$q1.Name = "2022-Q1"

# Match the page margins used throughout the rest of the workbook (the
# default for a brand-new sheet differs slightly).
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "166024"
$q1.Range("C2").Value = "中欧恒利三年定期开放混合"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "4.48"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "98.71"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "3.49"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.1564"
$q1.Range("H2").Value = 9

$q1.Range("A3").Value = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "005702"
$q1.Range("C3").Value = "恒生前海港股通高股息低波动指数"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "0.29"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "94.14"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "2.64"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.0077"
$q1.Range("H3").Value = 2

# Match the bold/centered/bordered header + index-column styling used by
# the other per-fund sheets.
$styleDonor.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)
$styleDonor.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Re-add "总计" (totals) at the end, with a new leading 2022-Q1 row.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$tot.Name = "总计"

$tot.PageSetup.LeftMargin = 54
$tot.PageSetup.RightMargin = 54
$tot.PageSetup.TopMargin = 72
$tot.PageSetup.BottomMargin = 72
$tot.PageSetup.HeaderMargin = 36
$tot.PageSetup.FooterMargin = 36
$tot.Outline.SummaryRow = 1
$tot.Outline.SummaryColumn = 1

$tot.Range("B1").Value = "日期"
$tot.Range("C1").Value = "持有数量(只)"
$tot.Range("D1").Value = "持有市值(亿元)"

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0.16

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2021-Q4"
$tot.Range("C3").Value = 3
$tot.Range("D3").Value = 0.06

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2021-Q3"
$tot.Range("C4").Value = 2
$tot.Range("D4").Value = 0.08

$tot.Range("A5").Value = 3
$tot.Range("B5").Value = "2020-Q4"
$tot.Range("C5").Value = 2
$tot.Range("D5").Value = 0.05

$styleDonor.Range("B1:D1").Copy()
$tot.Range("B1:D1").PasteSpecial($xlPasteFormats)
$styleDonor.Range("A2:A4").Copy()
$tot.Range("A2:A5").PasteSpecial($xlPasteFormats)
